# Update countries & provincias Spain
# 1) Refresh case numbers for a handful of countries.
# 2) Re-sort the whole country table by "Casos totales" (column B) descending,
#    keeping the original relative order for ties (stable sort).
# 3) Bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) apply the per-country value updates -------------------------------
$updates = @{
    "Iran"             = @{ D=6745;  E=11466 }
    "Alemania"         = @{ D=180;   E=17145 }
    "Estados Unidos"   = @{ B=14372; C=583; G=11; H=218 }
    "Suecia"           = @{ B=1639;  C=200;  E=1607; G=5; H=16 }
    "Brasil"           = @{ B=650;   C=10;   E=641 }
    "Chile"            = @{ B=434;   C=92;   D=6;   E=428; F=7;  G=0; H=0 }
    "Vietnam"          = @{ B=91;    C=6;    D=17;  E=74;  F=0;  G=0; H=0 }
    "Sri Lanka"        = @{ B=71;    C=11;   D=3;   E=68;  F=0;  G=0; H=0 }
    "Macao"            = @{ B=17;    C=0;    D=10;  E=7;   F=0;  G=0; H=0 }
    "Bolivia"          = @{ B=16;    C=1;    D=0;   E=16;  F=0;  G=0; H=0 }
}

$firstRow = 4
$lastRow = 186
$dataRange = $ws.Range("A" + $firstRow + ":H" + $lastRow)
$data = $dataRange.Value()

$rowCount = $data.GetLength(0)
$colNameToIndex = @{ "A"=1; "B"=2; "C"=3; "D"=4; "E"=5; "F"=6; "G"=7; "H"=8 }

for ($i = 1; $i -le $rowCount; $i++) {
    $country = $data[$i,1]
    if ($updates.ContainsKey($country)) {
        $fields = $updates[$country]
        foreach ($col in $fields.Keys) {
            $colIdx = $colNameToIndex[$col]
            $data[$i,$colIdx] = $fields[$col]
        }
    }
}

# ---- 2) stable sort by column B (Casos totales) descending ----------------
$rows = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $rowCount; $i++) {
    $rowVals = New-Object 'object[]' 8
    for ($j = 1; $j -le 8; $j++) {
        $rowVals[$j-1] = $data[$i,$j]
    }
    $entry = [PSCustomObject]@{ Total = [double]$data[$i,2]; Orig = $i; Row = $rowVals }
    $rows.Add($entry) | Out-Null
}

$sorted = $rows | Sort-Object -Property @{Expression="Total"; Descending=$true}

$outArr = New-Object 'object[,]' $rowCount,8
$r = 0
foreach ($entry in $sorted) {
    for ($j = 0; $j -lt 8; $j++) {
        $outArr[$r,$j] = $entry.Row[$j]
    }
    $r++
}

# Clear the country-name column first so the workbook's string pool drops every
# name that is about to be rewritten; re-entering them below (in final row
# order) makes the pool pick the names back up in that same first-use order,
# matching how a real re-sort/save in Excel rebuilds sharedStrings.xml.
$nameColRange = $ws.Range("A" + $firstRow + ":A" + $lastRow)
$nameColRange.Value = ""

for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item($firstRow + $i, 1).Value = $outArr[$i,0]
}

$dataRange.Value = $outArr

# ---- 3) bump the timestamp (written last so it sorts after the country
#         names in the rebuilt string pool, same as the source workbook) -----
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 14:46"
